$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2077922077922078
$ws.Range("C2").Value = 0.5411255411255411
$ws.Range("J2").Value = 0.02597402597402598
$ws.Range("P2").Value = 0.1558441558441558
$ws.Range("S2").Value = 0.06926406926406926
$ws.Range("C3").Value = 0.03053435114503817
$ws.Range("J3").Value = 0.03053435114503817
$ws.Range("P3").Value = 0.7786259541984732
$ws.Range("S3").Value = 0.1603053435114504
$ws.Range("J4").Value = 0.09090909090909091
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.2121212121212121
$ws.Range("B6").Value = 0.06542056074766354
$ws.Range("D6").Value = 0.009345794392523364
$ws.Range("F6").Value = 0.09345794392523364
$ws.Range("J6").Value = 0.2383177570093458
$ws.Range("O6").Value = 0.009345794392523364
$ws.Range("Q6").Value = 0.1448598130841121
$ws.Range("R6").Value = 0.07476635514018691
$ws.Range("S6").Value = 0.3644859813084112
$ws.Range("B7").Value = 0.1013824884792627
$ws.Range("D7").Value = 0.03686635944700461
$ws.Range("F7").Value = 0.07834101382488479
$ws.Range("J7").Value = 0.1013824884792627
$ws.Range("O7").Value = 0.02764976958525346
$ws.Range("Q7").Value = 0.184331797235023
$ws.Range("R7").Value = 0.05990783410138249
$ws.Range("S7").Value = 0.4101382488479263
$ws.Range("B8").Value = 0.08705882352941176
$ws.Range("D8").Value = 0.0188235294117647
$ws.Range("F8").Value = 0.04705882352941176
$ws.Range("J8").Value = 0.1082352941176471
$ws.Range("O8").Value = 0.02823529411764706
$ws.Range("Q8").Value = 0.1835294117647059
$ws.Range("R8").Value = 0.08941176470588236
$ws.Range("S8").Value = 0.4376470588235294
$ws.Range("B9").Value = 0.06349206349206349
$ws.Range("D9").Value = 0.01587301587301587
$ws.Range("F9").Value = 0.08994708994708994
$ws.Range("J9").Value = 0.09523809523809523
$ws.Range("O9").Value = 0.02645502645502645
$ws.Range("Q9").Value = 0.2116402116402116
$ws.Range("R9").Value = 0.06878306878306878
$ws.Range("S9").Value = 0.4285714285714285
$ws.Range("B10").Value = 0.08691588785046729
$ws.Range("D10").Value = 0.01214953271028037
$ws.Range("F10").Value = 0.07757009345794393
$ws.Range("J10").Value = 0.1364485981308411
$ws.Range("O10").Value = 0.02149532710280374
$ws.Range("Q10").Value = 0.202803738317757
$ws.Range("R10").Value = 0.08224299065420561
$ws.Range("S10").Value = 0.3803738317757009
$ws.Range("G11").Value = 0.1202749140893471
$ws.Range("J11").Value = 0.0584192439862543
$ws.Range("K11").Value = 0.1443298969072165
$ws.Range("L11").Value = 0.6597938144329897
$ws.Range("S11").Value = 0.01718213058419244
$ws.Range("G12").Value = 0.7563451776649747
$ws.Range("J12").Value = 0.1725888324873096
$ws.Range("K12").Value = 0.005076142131979695
$ws.Range("L12").Value = 0.03045685279187817
$ws.Range("S12").Value = 0.03553299492385787
$ws.Range("G13").Value = 0.6842105263157895
$ws.Range("J13").Value = 0.2807017543859649
$ws.Range("S13").Value = 0.03508771929824561
$ws.Range("F15").Value = 0.01104972375690608
$ws.Range("H15").Value = 0.1104972375690608
$ws.Range("I15").Value = 0.06629834254143646
$ws.Range("J15").Value = 0.3204419889502763
$ws.Range("K15").Value = 0.08839779005524862
$ws.Range("M15").Value = 0.01657458563535912
$ws.Range("O15").Value = 0.09392265193370165
$ws.Range("S15").Value = 0.292817679558011
$ws.Range("F16").Value = 0.04487179487179487
$ws.Range("H16").Value = 0.2115384615384615
$ws.Range("I16").Value = 0.08333333333333333
$ws.Range("J16").Value = 0.3782051282051282
$ws.Range("K16").Value = 0.1153846153846154
$ws.Range("M16").Value = 0.00641025641025641
$ws.Range("O16").Value = 0.02564102564102564
$ws.Range("S16").Value = 0.1346153846153846
$ws.Range("F17").Value = 0.01485148514851485
$ws.Range("H17").Value = 0.2128712871287129
$ws.Range("I17").Value = 0.09653465346534654
$ws.Range("J17").Value = 0.3737623762376238
$ws.Range("K17").Value = 0.1064356435643564
$ws.Range("M17").Value = 0.0297029702970297
$ws.Range("O17").Value = 0.04207920792079208
$ws.Range("S17").Value = 0.1237623762376238
$ws.Range("F18").Value = 0.01197604790419162
$ws.Range("H18").Value = 0.09580838323353294
$ws.Range("I18").Value = 0.1197604790419162
$ws.Range("J18").Value = 0.4670658682634731
$ws.Range("K18").Value = 0.1497005988023952
$ws.Range("M18").Value = 0.04191616766467066
$ws.Range("O18").Value = 0.04790419161676647
$ws.Range("S18").Value = 0.0658682634730539
$ws.Range("F19").Value = 0.01595298068849706
$ws.Range("H19").Value = 0.2250209907640638
$ws.Range("I19").Value = 0.0873215785054576
$ws.Range("J19").Value = 0.3249370277078086
$ws.Range("K19").Value = 0.1192275398824517
$ws.Range("M19").Value = 0.02854743912678422
$ws.Range("O19").Value = 0.05625524769101595
$ws.Range("S19").Value = 0.1427371956339211
